# Rename worksheet tabs to new spatial-unit summary IDs
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ11373791",
    "summ12375527",
    "summ13501561",
    "summ14519540",
    "summ15577260",
    "summ16891102",
    "summ17888260",
    "summ19092725",
    "summ20212789"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
